$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "CO2 In (Predicted DNN)"
$ws.Range("J1").Value = "Temperature In (Predicted DNN)"
$ws.Range("N1").Value = "RH In (Predicted DNN)"
$ws.Range("R1").Value = "PAR In (Predicted DNN)"
$ws.Range("V1").Value = "Leaf Temp (Predicted DNN)"
